# TC10_Canine_Filter_Breed-BassHnd.xlsx -- "startup" sheet edit
#
# The cell B2 (style s="2", wrap-text) previously held no value; it now
# carries the Neo4j/Cypher query used to pull the "Basset Hound" filtered
# data set. Setting the value also grows the (wrap-text) row to fit the
# long query text, and the sheet's view / selection is left parked on B2
# (matching the saved cursor position recorded in the workbook).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

$query = "MATCH (s:study) WITH COLLECT(DISTINCT(s.clinical_study_designation)) AS all_studies MATCH (d:demographic) WITH COLLECT(DISTINCT(d.breed)) AS all_breeds, COLLECT(DISTINCT(d.sex)) AS all_sexes, all_studies MATCH (d:diagnosis) WITH COLLECT(DISTINCT(d.disease_term)) AS all_diseases, all_breeds, all_sexes, all_studies MATCH (p:program)<-[*]-(s:study)<-[*]-(c:case)<--(demo:demographic), (c)<--(diag:diagnosis) WHERE demo.breed IN ['Basset Hound'] WITH DISTINCT c AS c, p, s, demo, diag RETURN coalesce(c.case_id,'') AS ``Case ID`` , coalesce(s.clinical_study_designation,'') AS ``Study Code`` , coalesce(s.clinical_study_type,'') AS  ``Study Type``, coalesce(demo.breed,'') AS Breed , coalesce(diag.disease_term,'') AS Diagnosis , coalesce(diag.stage_of_disease,'') AS ``Stage of Disease`` ,  coalesce(demo.patient_age_at_enrollment,'') AS Age , coalesce(demo.sex,'') AS Sex , coalesce(demo.neutered_indicator,'') AS  ``Neutered Status``"

# Write the Cypher query text into B2 (adds a new shared string).
$ws.Range("B2").Value = $query

# B2 wraps text; grow row 2 so the whole query is visible (~188.5pt).
$ws.Rows.Item(2).RowHeight = 188.5

# Park the selection/view on B2.
$ws.Activate()
$ws.Range("B2").Select()
